$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A97").Value = "2025-04-29 16:38:35"
$ws.Range("B97").Value = 252
